# Task1 (Simple Sequence Analysis) completion.
#
# The workbook has two worksheets:
#   "forward" (sheet1) - DNA/DNA2/mRNA/tRNA/Protein table read 5'->3'
#   "reverse" (sheet2) - the same table for the complementary strand
# Both contain a codon-by-codon translation exercise where every third
# row's column E holds the translated amino acid (3-letter code / STOP).
# This edit fills in / corrects the amino-acid answers for both tables,
# and adds a short note next to the start codon on the "forward" sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "forward" (sheet1): fill in / correct the Protein column
# ------------------------------------------------------------------
$wsFwd = $wb.Worksheets.Item(1)

$wsFwd.Range("E10").Value = "Met"
$wsFwd.Range("F10").Value = "Begynner på Met!!"
$wsFwd.Range("E13").Value = "Ser"
$wsFwd.Range("E16").Value = "Gly"
$wsFwd.Range("E19").Value = "His"
$wsFwd.Range("E22").Value = "Leu"
$wsFwd.Range("E25").Value = "Pro"
$wsFwd.Range("E28").Value = "Arg"
$wsFwd.Range("E31").Value = "Thr"

[void]$wsFwd.Range("G11").Select()

# ------------------------------------------------------------------
# Sheet "reverse" (sheet2): move / correct the Protein column answers
# ------------------------------------------------------------------
$wsRev = $wb.Worksheets.Item(2)

$wsRev.Range("E5").Value = ""
$wsRev.Range("E7").Value = "Arg"
$wsRev.Range("E8").Value = ""
$wsRev.Range("E10").Value = "Tyr"
$wsRev.Range("E11").Value = ""
$wsRev.Range("E13").Value = "STOP"
$wsRev.Range("E14").Value = ""
$wsRev.Range("E16").Value = "Ala"
$wsRev.Range("E17").Value = ""
$wsRev.Range("E19").Value = "Val"
$wsRev.Range("E20").Value = ""
$wsRev.Range("E22").Value = "Gln"
$wsRev.Range("E23").Value = ""
$wsRev.Range("E25").Value = "Arg"
$wsRev.Range("E26").Value = ""
$wsRev.Range("E28").Value = "Pro"
$wsRev.Range("E29").Value = ""
$wsRev.Range("E31").Value = "Ser"

[void]$wsRev.Range("F10").Select()

[void]$wsFwd.Activate()
